$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that, if written as plain numeric-looking text (single decimal point),
# would be auto-converted to numbers by Excel. Force Text format first so the literal
# string is preserved, then restore the default "Normal" style so no stray formatting
# is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}


# Row 2
$ws.Range("D2").Value = "60.577.76"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3
$ws.Range("D3").Value = "2.625.58"
$ws.Range("E3").Value = "  +1.12%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "581.27"
$ws.Range("E5").Value = "  +2.30%  "

# Row 6
Set-TextValue $ws.Range("D6") "144.68"
$ws.Range("E6").Value = "  +1.88%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 10
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("E11").Value = "  +2.03%  "

# Row 12
$ws.Range("E12").Value = "  +3.38%  "

# Row 13
$ws.Range("D13").Value = "3.103.95"
$ws.Range("E13").Value = "  +1.50%  "

# Row 14
Set-TextValue $ws.Range("D14") "26.07"
$ws.Range("E14").Value = "  +11.35%  "

# Row 15
$ws.Range("D15").Value = "60.568.79"

# Row 16
$ws.Range("E16").Value = "  +1.23%  "

# Row 17
$ws.Range("D17").Value = "2.632.66"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.53"
$ws.Range("E18").Value = "  +2.40%  "

# Row 19
$ws.Range("E19").Value = "  +1.15%  "

# Row 20
Set-TextValue $ws.Range("D20") "347.99"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.89"
$ws.Range("E21").Value = "  -1.50%  "

# Row 22
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
Set-TextValue $ws.Range("D24") "63.91"
$ws.Range("E24").Value = "  +1.24%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.999"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26
$ws.Range("E26").Value = "  +1.74%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.12"
$ws.Range("E27").Value = "  +5.84%  "

# Row 28
$ws.Range("E28").Value = "  +11.34%  "

# Row 29
$ws.Range("E29").Value = "  +1.60%  "

# Row 30
Set-TextValue $ws.Range("D30") "6.67"
$ws.Range("E30").Value = "  +5.41%  "

# Row 31
Set-TextValue $ws.Range("D31") "168.52"
$ws.Range("E31").Value = "  +4.72%  "

# Row 32
$ws.Range("E32").Value = "  +0.11%  "

# Row 33
$ws.Range("E33").Value = "  +0.87%  "

# Row 34
$ws.Range("E34").Value = "  +10.05%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.42"
$ws.Range("E35").Value = "  +4.92%  "

# Row 36
$ws.Range("E36").Value = "  +7.97%  "

# Row 37
$ws.Range("E37").Value = "  +2.50%  "

# Row 38
Set-TextValue $ws.Range("D38") "331.96"
$ws.Range("E38").Value = "  +12.91%  "

# Row 39
$ws.Range("E39").Value = "  +4.95%  "

# Row 40
Set-TextValue $ws.Range("D40") "38.46"
$ws.Range("E40").Value = "  +1.70%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.874"
$ws.Range("E41").Value = "  +3.04%  "

# Row 42
Set-TextValue $ws.Range("D42") "5.15"
$ws.Range("E42").Value = "  +6.93%  "

# Row 43
Set-TextValue $ws.Range("D43") "20.69"
$ws.Range("E43").Value = "  +4.11%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D44") "133.11"
$ws.Range("E44").Value = "  -3.52%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D45") "0.0997"
$ws.Range("E45").Value = "  +1.69%  "

# Row 46
Set-TextValue $ws.Range("D46") "20.09"
$ws.Range("E46").Value = "  +2.08%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.999"
$ws.Range("E47").Value = "  +0.29%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0557"
$ws.Range("E48").Value = "  +2.05%  "

# Row 49
$ws.Range("E49").Value = "  +0.31%  "

# Row 50
$ws.Range("E50").Value = "  +1.98%  "

# Row 51
$ws.Range("E51").Value = "  +0.55%  "
